# Warehouse_Progress milestone data fix
# - Landlord Refurb Design Approved: shift target/actual dates back one month
# - "Landlord Refurb In Progress" -> "Warehouse Refurb", with updated start
#   date, progress %, and note text (now correctly "in progress" at 15%)
# - Internal Fit-out Planning: target date pushed out (after the refurb)
# - Racking & Storage Installed: target date pulled forward

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warehouse_Progress")

# These "date" columns (Target Date / Actual Date) are stored as plain text
# (e.g. "2024-11-01"), not real Excel dates. Assigning a string that looks
# like a date through .Value would normally get auto-converted into a date
# serial by Excel's General-format type inference, so force the cell to
# Text first, then clear the format override afterwards so no stray
# NumberFormat is left behind on the cell.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 3: Landlord Refurb Design Approved
Set-TextValue $ws.Range("B3") "2024-11-01"
Set-TextValue $ws.Range("E3") "2024-10-28"

# Row 4: Landlord Refurb In Progress -> Warehouse Refurb
$ws.Range("A4").Value = "Warehouse Refurb"
Set-TextValue $ws.Range("B4") "2025-08-01"
$ws.Range("C4").Value = 15
$ws.Range("F4").Value = "Started 4 Nov 2024 - landlord managing refurbishment (4 weeks in)"

# Row 5: Internal Fit-out Planning
Set-TextValue $ws.Range("B5") "2025-09-01"
$ws.Range("F5").Value = "Racking layout design - begins as refurb nears completion"

# Row 6: Racking & Storage Installed
Set-TextValue $ws.Range("B6") "2025-11-15"
